$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value() = 3310
$ws.Range("I12").Value() = 3310
$ws.Range("J12").Value() = 0
$ws.Range("K12").Value() = 3310
$ws.Range("L12").Value() = 0
$ws.Range("M12").Value() = -3140
$ws.Range("N12").ClearContents()

$ws.Range("H19").Value() = 4667.25
$ws.Range("J19").Value() = 1447.5
$ws.Range("L19").Value() = 1447.5
$ws.Range("N19").Value() = -1797.5

$ws.Range("H42").Value() = 932.625
$ws.Range("I42").Value() = 72.666664
$ws.Range("J42").Value() = 3512.5
$ws.Range("K42").Value() = 217.999992
$ws.Range("L42").Value() = 10537.5
$ws.Range("M42").Value() = 12.00000800000001
$ws.Range("N42").Value() = -10997.5

$ws.Range("H58").Value() = 955
$ws.Range("I58").Value() = 513.25
$ws.Range("K58").Value() = 1539.75
$ws.Range("M58").Value() = -1389.75

$ws.Range("H62").Value() = 8477.857
$ws.Range("I62").Value() = 7136.8
$ws.Range("K62").Value() = 7136.8
$ws.Range("M62").Value() = -6512.8

$ws.Range("H65").Value() = 8477.857
$ws.Range("I65").Value() = 7136.8
$ws.Range("K65").Value() = 35684
$ws.Range("M65").Value() = -32564

$ws.Range("H76").Value() = 6699.4
$ws.Range("I76").Value() = 1997.5
$ws.Range("K76").Value() = 1997.5
$ws.Range("M76").Value() = -1682.5

$ws.Range("H79").Value() = 6699.4
$ws.Range("I79").Value() = 1997.5
$ws.Range("K79").Value() = 1997.5
$ws.Range("M79").Value() = -905.5

$ws.Range("H131").Value() = 4929.3335
$ws.Range("J131").Value() = 7262.375
$ws.Range("L131").Value() = 21787.125
$ws.Range("N131").Value() = -31867.125

$ws.Range("H138").Value() = 4130.523
$ws.Range("J138").Value() = 4235
$ws.Range("L138").Value() = 12705
$ws.Range("N138").Value() = -22985

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value() = 1491.8928
$ws.Range("I32").Value() = 1518.7
$ws.Range("J32").Value() = 1268.5
$ws.Range("K32").Value() = 1518.7
$ws.Range("L32").Value() = 1268.5
$ws.Range("M32").Value() = -1231.7
$ws.Range("N32").Value() = -1842.5

$ws.Range("H45").Value() = 7214.933
$ws.Range("I45").Value() = 5767.4165
$ws.Range("K45").Value() = 5767.4165
$ws.Range("M45").Value() = -5390.4165

$ws.Range("H61").Value() = 6399.1055
$ws.Range("I61").Value() = 5532.3887
$ws.Range("J61").Value() = 22000
$ws.Range("K61").Value() = 5532.3887
$ws.Range("L61").Value() = 22000
$ws.Range("M61").Value() = -5320.3887
$ws.Range("N61").Value() = -22424

$ws.Range("H122").Value() = 4486.35
$ws.Range("I122").Value() = 3911.7273
$ws.Range("K122").Value() = 11735.1819
$ws.Range("M122").Value() = -9285.1819

$ws.Range("H136").Value() = 6399.1055
$ws.Range("I136").Value() = 5532.3887
$ws.Range("J136").Value() = 22000
$ws.Range("K136").Value() = 16597.1661
$ws.Range("L136").Value() = 66000
$ws.Range("M136").Value() = -14047.1661
$ws.Range("N136").Value() = -71100

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value() = 17
$ws.Range("I4").Value() = 17
$ws.Range("J4").Value() = 0
$ws.Range("K4").Value() = 17
$ws.Range("L4").Value() = 0
$ws.Range("M4").Value() = 95
$ws.Range("N4").ClearContents()

$ws.Range("H31").Value() = 17533.12
$ws.Range("I31").Value() = 2685.0454
$ws.Range("J31").Value() = 23696.473
$ws.Range("K31").Value() = 2685.0454
$ws.Range("L31").Value() = 23696.473
$ws.Range("M31").Value() = -2390.0454
$ws.Range("N31").Value() = -24286.473

$ws.Range("H34").Value() = 17533.12
$ws.Range("I34").Value() = 2685.0454
$ws.Range("J34").Value() = 23696.473
$ws.Range("K34").Value() = 2685.0454
$ws.Range("L34").Value() = 23696.473
$ws.Range("M34").Value() = -2483.0454
$ws.Range("N34").Value() = -24100.473

$ws.Range("H86").Value() = 9023.299999999999
$ws.Range("I86").Value() = 4587.25
$ws.Range("J86").Value() = 11980.667
$ws.Range("K86").Value() = 4587.25
$ws.Range("L86").Value() = 11980.667
$ws.Range("M86").Value() = -3464.25
$ws.Range("N86").Value() = -14226.667

$ws.Range("H89").Value() = 9023.299999999999
$ws.Range("I89").Value() = 4587.25
$ws.Range("J89").Value() = 11980.667
$ws.Range("K89").Value() = 22936.25
$ws.Range("L89").Value() = 59903.335
$ws.Range("M89").Value() = -17320.25
$ws.Range("N89").Value() = -71135.33499999999

$ws.Range("H122").Value() = 4622.8486
$ws.Range("I122").Value() = 2413.5
$ws.Range("J122").Value() = 8021.846
$ws.Range("K122").Value() = 7240.5
$ws.Range("L122").Value() = 24065.538
$ws.Range("M122").Value() = -4790.5
$ws.Range("N122").Value() = -28965.538

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H9").Value() = 1828333.5
$ws.Range("J9").Value() = 5000
$ws.Range("L9").Value() = 15000
$ws.Range("N9").Value() = -15448

$ws.Range("H41").Value() = 148.33333
$ws.Range("I41").Value() = 148.33333
$ws.Range("J41").Value() = 0
$ws.Range("K41").Value() = 444.99999
$ws.Range("L41").Value() = 0
$ws.Range("M41").Value() = -106.99999
$ws.Range("N41").ClearContents()

$ws.Range("H92").Value() = 258.7143
$ws.Range("J92").Value() = 268.57693
$ws.Range("L92").Value() = 805.7307900000001
$ws.Range("N92").Value() = -3301.73079

$ws.Range("H107").Value() = 1443.8
$ws.Range("J107").Value() = 2472.5557
$ws.Range("L107").Value() = 7417.6671
$ws.Range("N107").Value() = -11257.6671

$ws.Range("H129").Value() = 5210863
$ws.Range("I129").Value() = 614.1
$ws.Range("J129").Value() = 13894611
$ws.Range("K129").Value() = 1842.3
$ws.Range("L129").Value() = 41683833
$ws.Range("M129").Value() = 3157.7
$ws.Range("N129").Value() = -41693833

$ws.Range("H131").Value() = 6483861.5
$ws.Range("J131").Value() = 14960341
$ws.Range("L131").Value() = 44881023
$ws.Range("N131").Value() = -44891103

$ws.Range("H132").Value() = 3139.121
$ws.Range("I132").Value() = 3000.5557
$ws.Range("J132").Value() = 3305.4
$ws.Range("K132").Value() = 27005.0013
$ws.Range("L132").Value() = 29748.6
$ws.Range("M132").Value() = -24475.0013
$ws.Range("N132").Value() = -34808.60000000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value() = 0
$ws.Range("I5").Value() = 0
$ws.Range("J5").Value() = 0
$ws.Range("K5").Value() = 0
$ws.Range("L5").Value() = 0
$ws.Range("M5").ClearContents()
$ws.Range("N5").ClearContents()

$ws.Range("H53").Value() = 43999
$ws.Range("I53").Value() = 44750
$ws.Range("K53").Value() = 44750
$ws.Range("M53").Value() = -44119

$ws.Range("H80").Value() = 4295.84
$ws.Range("I80").Value() = 2954.9
$ws.Range("K80").Value() = 2954.9
$ws.Range("M80").Value() = -1956.9

$ws.Range("H83").Value() = 4295.84
$ws.Range("I83").Value() = 2954.9
$ws.Range("K83").Value() = 14774.5
$ws.Range("M83").Value() = -9782.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value() = 2078
$ws.Range("I16").Value() = 2086.3
$ws.Range("K16").Value() = 2086.3
$ws.Range("M16").Value() = -1916.3

$ws.Range("H40").Value() = 16321
$ws.Range("I40").Value() = 15921.0625
$ws.Range("K40").Value() = 15921.0625
$ws.Range("M40").Value() = -15785.0625

$ws.Range("H46").Value() = 4527.2856
$ws.Range("I46").Value() = 0
$ws.Range("J46").Value() = 4527.2856
$ws.Range("K46").Value() = 0
$ws.Range("L46").Value() = 4527.2856
$ws.Range("N46").Value() = -4903.2856
$ws.Range("M46").ClearContents()

$ws.Range("H55").Value() = 2632235
$ws.Range("I55").Value() = 5000197.5
$ws.Range("J55").Value() = 1165.5555
$ws.Range("K55").Value() = 5000197.5
$ws.Range("L55").Value() = 1165.5555
$ws.Range("M55").Value() = -5000024.5
$ws.Range("N55").Value() = -1511.5555

$ws.Range("H61").Value() = 5916
$ws.Range("I61").Value() = 3874
$ws.Range("K61").Value() = 3874
$ws.Range("M61").Value() = -3672

$ws.Range("H113").Value() = 5916
$ws.Range("I113").Value() = 3874
$ws.Range("K113").Value() = 3874
$ws.Range("M113").Value() = -1704

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H26").Value() = 2798.3
$ws.Range("I26").Value() = 2887
$ws.Range("J26").Value() = 2000
$ws.Range("K26").Value() = 2887
$ws.Range("L26").Value() = 2000
$ws.Range("M26").Value() = -2594
$ws.Range("N26").Value() = -2586

$ws.Range("H81").Value() = 12768.333
$ws.Range("I81").Value() = 1001
$ws.Range("K81").Value() = 2002
$ws.Range("M81").Value() = -941

$ws.Range("H84").Value() = 12768.333
$ws.Range("I84").Value() = 1001
$ws.Range("K84").Value() = 10010
$ws.Range("M84").Value() = -4706
